# Applies the v1.5 bitacora update:
#  - Log sheet: 3 new rows (30-32) describing tipo de error / duplicados / filtro
#  - Resumen sheet: 4 new rows (31-34) describing the same features
#  - Versiones sheet: 1 new row (7) for version 1.5

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Log": add rows 30-32
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$logRows = @(
    @("27/02/2025", "18:30", "Tipo de error y detección de potencial duplicado", "En la solapa Errores: columna Tipo de error (Inconsistencia entre Categoria/Cuenta/Descripcion o Potencial registro duplicado). Detección de duplicados por misma fecha, monto, tipo_movimiento y descripción similar. Para duplicados: icono Ver que abre modal comparando ambos registros; opciones Excluir de cálculos (anular) o Eliminar registro. Export Excel incluye tipo_error.", "Diagnostico"),
    @("27/02/2025", "18:40", "Filtro por tipo de error en solapa Errores", "Selector ""Tipo de error"" en la barra de la solapa Errores: Todos, Inconsistencia (categoría/cuenta/descripción), Potencial registro duplicado. La tabla y la exportación a Excel respetan el filtro seleccionado.", "Diagnostico"),
    @("27/02/2025", "18:50", "Duplicados: cliente igual e id_origen en comparación", "Solo se marca potencial duplicado si además de fecha, monto, tipo y descripción similar el campo cliente es igual; si cliente es distinto no se marca. En el modal de comparación (Este registro / Posible duplicado) se incluye id_origen y Cliente.", "Diagnostico")
)

$r = 30
foreach ($row in $logRows) {
    $log.Cells.Item($r, 1).Value = $row[0]
    $log.Cells.Item($r, 2).Value = $row[1]
    $log.Cells.Item($r, 3).Value = $row[2]
    $log.Cells.Item($r, 4).Value = $row[3]
    $log.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "Resumen": add rows 31-34
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

$resumenRows = @(
    @("Tipo de error en Errores", "Tabla de errores muestra columna Tipo de error: Inconsistencia entre Categoria, Cuenta Contable y Descripcion; o Potencial registro duplicado. Export a Excel incluye tipo_error."),
    @("Detección de potencial duplicado", "Registros con misma fecha, monto, tipo_movimiento y descripción similar se marcan como potencial duplicado. Icono Ver abre modal con comparación Este registro / Posible duplicado; acciones: Excluir de cálculos (anular) o Eliminar registro."),
    @("Filtro por tipo de error", "En la solapa Errores, selector para filtrar por tipo: Todos, Inconsistencia (categoría/cuenta/descripción), Potencial registro duplicado. La exportación a Excel exporta solo los registros visibles según el filtro."),
    @("Duplicados: condición cliente", "Dos registros son potencial duplicado solo si coinciden en fecha, monto, tipo_movimiento, descripción similar y además cliente es igual; si cliente es distinto no se marcan como duplicado. Modal de comparación muestra id_origen y Cliente.")
)

$r = 31
foreach ($row in $resumenRows) {
    $resumen.Cells.Item($r, 1).Value = $row[0]
    $resumen.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "Versiones": add row 7
# ---------------------------------------------------------------------------
$versiones = $wb.Worksheets.Item("Versiones")

$versiones.Cells.Item(7, 1).Value = "1.5"
$versiones.Cells.Item(7, 2).Value = "27/02/2025"
$versiones.Cells.Item(7, 3).Value = "Errores: tipo de error, detección duplicados (cliente igual), filtro por tipo, modal comparación con id_origen; timeout carga y fechaStr para fechas"
